$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps its scraped text formatting (e.g. "299.01",
# "42.203.06") instead of being auto-converted to a number by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "42.203.06"
$ws.Range("E2").Value = "  -1.42%  "
$ws.Range("D3").Value = "2.270.54"
$ws.Range("E3").Value = "  -2.29%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "299.01"
$ws.Range("E5").Value = "  -1.96%  "
$ws.Range("D6").Value = "95.76"
$ws.Range("E6").Value = "  -4.51%  "
$ws.Range("D7").Value = "0.496"
$ws.Range("E7").Value = "  -2.40%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").Value = "0.494"
$ws.Range("E9").Value = "  -1.89%  "
$ws.Range("D10").Value = "33.43"
$ws.Range("E10").Value = "  -2.78%  "
$ws.Range("D11").Value = "0.0791"
$ws.Range("E11").Value = "  -0.12%  "
$ws.Range("D12").Value = "48.18"
$ws.Range("E12").Value = "  -7.57%  "
$ws.Range("E13").Value = "  +0.01%  "
$ws.Range("D14").Value = "6.67"
$ws.Range("E14").Value = "  -1.01%  "
$ws.Range("D15").Value = "2.624.30"
$ws.Range("E15").Value = "  -2.52%  "
$ws.Range("D16").Value = "15.56"
$ws.Range("E16").Value = "  -0.65%  "
$ws.Range("D17").Value = "2.274.69"
$ws.Range("E17").Value = "  -2.38%  "
$ws.Range("D18").Value = "0.784"
$ws.Range("E18").Value = "  -4.40%  "
$ws.Range("D19").Value = "42.145.04"
$ws.Range("E19").Value = "  -1.35%  "
$ws.Range("D20").Value = "11.73"
$ws.Range("E20").Value = "  +1.53%  "
$ws.Range("D21").Value = "0.0₃0892"
$ws.Range("E21").Value = "  -1.08%  "
$ws.Range("D22").Value = "5.99"
$ws.Range("E22").Value = "  -2.38%  "
$ws.Range("D23").Value = "66.65"
$ws.Range("E23").Value = "  -3.46%  "
$ws.Range("D24").Value = "235.04"
$ws.Range("E24").Value = "  -0.07%  "
$ws.Range("D25").Value = "1.96"
$ws.Range("E25").Value = "  -0.84%  "
$ws.Range("E26").Value = "  +0.30%  "
$ws.Range("D27").Value = "2.45"
$ws.Range("E27").Value = "  -3.31%  "
$ws.Range("D28").Value = "24.00"
$ws.Range("E28").Value = "  -4.98%  "
$ws.Range("D29").Value = "2.29"
$ws.Range("E29").Value = "  +5.02%  "
$ws.Range("D30").Value = "168.29"
$ws.Range("E30").Value = "  +4.75%  "
$ws.Range("D31").Value = "33.97"
$ws.Range("E31").Value = "  -2.11%  "
$ws.Range("D32").Value = "9.15"
$ws.Range("E32").Value = "  -0.47%  "
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  -0.10%  "
$ws.Range("D34").Value = "4.91"
$ws.Range("E34").Value = "  -2.68%  "
$ws.Range("D35").Value = "4.52"
$ws.Range("E35").Value = "  -1.13%  "
$ws.Range("E36").Value = "  -2.12%  "
$ws.Range("E37").Value = "  -4.86%  "
$ws.Range("E38").Value = "  -4.04%  "
$ws.Range("D39").Value = "2.78"
$ws.Range("E39").Value = "  -3.11%  "
$ws.Range("D40").Value = "0.0987"
$ws.Range("E40").Value = "  -2.27%  "
$ws.Range("E41").Value = "  -2.26%  "
$ws.Range("E42").Value = "  -5.07%  "
$ws.Range("D43").Value = "2.44"
$ws.Range("E43").Value = "  -5.06%  "
$ws.Range("D44").Value = "1.960.47"
$ws.Range("E44").Value = "  -2.21%  "
$ws.Range("E45").Value = "  -1.16%  "
$ws.Range("D46").Value = "17.42"
$ws.Range("E46").Value = "  -6.03%  "
$ws.Range("E47").Value = "  -5.82%  "
$ws.Range("D48").Value = "2.79"
$ws.Range("E48").Value = "  -2.76%  "
$ws.Range("D49").Value = "2.495.79"
$ws.Range("E49").Value = "  -2.31%  "
$ws.Range("D50").Value = "52.39"
$ws.Range("E50").Value = "  -5.13%  "
$ws.Range("E51").Value = "  -2.58%  "
